# Add data for 2022-04-26
# - Rename sheet / update header label from "April 17" to "April 18"
# - Update carjacking counts for several neighborhoods to reflect the
#   newly added day of data (2022-04-26 reporting window)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet itself (updates <sheet name="..."/> in workbook.xml)
$ws.Name = "Through 2022-04-18"

# Update the column header label (shared string used by cell B1)
$ws.Range("B1").Value = "April 2022 (through April 18)"

# --- Updated existing counts -------------------------------------------------
$ws.Range("J2").Value = 7     # Austin
$ws.Range("V2").Value = 5     # Austin
$ws.Range("B3").Value = 7     # Englewood
$ws.Range("N5").Value = 4     # Garfield Park
$ws.Range("B8").Value = 5     # Chicago Lawn
$ws.Range("R17").Value = 2    # Belmont Cragin
$ws.Range("J26").Value = 3    # South Shore
$ws.Range("F27").Value = 2    # Uptown
$ws.Range("B32").Value = 3    # Roseland

# --- Newly added counts (previously empty cells) -----------------------------
$ws.Range("F9").Value = 1     # Loop
$ws.Range("R26").Value = 1    # South Shore
$ws.Range("V49").Value = 1    # Mckinley Park
$ws.Range("Z50").Value = 1    # Grand Crossing
$ws.Range("Z57").Value = 1    # Albany Park
$ws.Range("F62").Value = 1    # Beverly
$ws.Range("Z73").Value = 1    # Jefferson Park
$ws.Range("V94").Value = 1    # West Ridge
